$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "April16"
$ws.Range("B2").Value = "reactionTime"
$ws.Range("C2").Value = [double]"574.5"
$ws.Range("D2").Value = [double]"4.731507720099267e-08"
$ws.Range("E2").Value = $true
$ws.Range("F2").Value = [double]"0.6059670781893004"
$ws.Range("G2").Value = [double]"4.649771442666329e-08"
$ws.Range("H2").Value = "Dunn"
$ws.Range("I2").Value = [double]"-1.046296296296299"
$ws.Range("J2").Value = "mannwhitneyu"

$ws.Range("A3").Value = "April16"
$ws.Range("B3").Value = "peakTime"
$ws.Range("C3").Value = [double]"1702"
$ws.Range("D3").Value = [double]"0.122235035799038"
$ws.Range("E3").Value = $false
$ws.Range("F3").Value = [double]"-0.1673525377229081"
$ws.Range("G3").Value = [double]"0.121469863914346"
$ws.Range("H3").Value = "Dunn"
$ws.Range("I3").Value = [double]"0.7339506172839521"
$ws.Range("J3").Value = "mannwhitneyu"

$ws.Range("A4").Value = "April16"
$ws.Range("B4").Value = "difference"
$ws.Range("C4").Value = [double]"2104"
$ws.Range("D4").Value = [double]"6.944883517617799e-05"
$ws.Range("E4").Value = $true
$ws.Range("F4").Value = [double]"-0.4430727023319616"
$ws.Range("G4").Value = [double]"6.855448672896294e-05"
$ws.Range("H4").Value = "Dunn"
$ws.Range("I4").Value = [double]"1.780246913580246"
$ws.Range("J4").Value = "mannwhitneyu"

$ws.Range("A5").Value = "April16"
$ws.Range("B5").Value = "peakValue"
$ws.Range("C5").Value = [double]"2881"
$ws.Range("D5").Value = [double]"2.319259952975573e-18"
$ws.Range("E5").Value = $true
$ws.Range("F5").Value = [double]"-0.9759945130315502"
$ws.Range("G5").Value = [double]"2.257024464666691e-18"
$ws.Range("H5").Value = "Dunn"
$ws.Range("I5").Value = [double]"61.07962962962966"
$ws.Range("J5").Value = "mannwhitneyu"

$ws.Range("A6").Value = "April16"
$ws.Range("B6").Value = "RMS"
$ws.Range("C6").Value = [double]"2869"
$ws.Range("D6").Value = [double]"4.44515293518166e-18"
$ws.Range("E6").Value = $true
$ws.Range("F6").Value = [double]"-0.9677640603566529"
$ws.Range("G6").Value = [double]"4.32684049725133e-18"
$ws.Range("H6").Value = "Dunn"
$ws.Range("I6").Value = [double]"16.35589506172839"
$ws.Range("J6").Value = "mannwhitneyu"

$ws.Range("A7").Value = "April16"
$ws.Range("B7").Value = "tau"
$ws.Range("C7").Value = [double]"554"
$ws.Range("D7").Value = [double]"2.831991584422178e-08"
$ws.Range("E7").Value = $true
$ws.Range("F7").Value = [double]"0.6200274348422496"
$ws.Range("G7").Value = [double]"2.782633665191067e-08"
$ws.Range("H7").Value = "Dunn"
$ws.Range("I7").Value = [double]"-21.61515123456792"
$ws.Range("J7").Value = "mannwhitneyu"

$ws.Range("A8").Value = "April16"
$ws.Range("B8").Value = "AUC"
$ws.Range("C8").Value = [double]"2802"
$ws.Range("D8").Value = [double]"1.518005587749135e-16"
$ws.Range("E8").Value = $true
$ws.Range("F8").Value = [double]"-0.9218106995884774"
$ws.Range("G8").Value = [double]"1.479448290688163e-16"
$ws.Range("H8").Value = "Dunn"
$ws.Range("I8").Value = [double]"5809.321388888889"
$ws.Range("J8").Value = "mannwhitneyu"

$ws.Range("A9").Value = "June26"
$ws.Range("B9").Value = "reactionTime"
$ws.Range("C9").Value = [double]"1575.5"
$ws.Range("D9").Value = [double]"0.4525480764086293"
$ws.Range("E9").Value = $false
$ws.Range("F9").Value = [double]"-0.08058984910836764"
$ws.Range("G9").Value = [double]"0.4506187202122447"
$ws.Range("H9").Value = "Dunn"
$ws.Range("I9").Value = [double]"0.2765432098765412"
$ws.Range("J9").Value = "mannwhitneyu"

$ws.Range("A10").Value = "June26"
$ws.Range("B10").Value = "peakTime"
$ws.Range("C10").Value = [double]"1183.5"
$ws.Range("D10").Value = [double]"0.08026984395540572"
$ws.Range("E10").Value = $false
$ws.Range("F10").Value = [double]"0.1882716049382716"
$ws.Range("G10").Value = [double]"0.07971977066430727"
$ws.Range("H10").Value = "Dunn"
$ws.Range("I10").Value = [double]"-0.8339506172839535"
$ws.Range("J10").Value = "mannwhitneyu"

$ws.Range("A11").Value = "June26"
$ws.Range("B11").Value = "difference"
$ws.Range("C11").Value = [double]"1277.5"
$ws.Range("D11").Value = [double]"0.2659278918781623"
$ws.Range("E11").Value = $false
$ws.Range("F11").Value = [double]"0.1237997256515775"
$ws.Range("G11").Value = [double]"0.2646022141905987"
$ws.Range("H11").Value = "Dunn"
$ws.Range("I11").Value = [double]"-1.110493827160496"
$ws.Range("J11").Value = "mannwhitneyu"

$ws.Range("A12").Value = "June26"
$ws.Range("B12").Value = "peakValue"
$ws.Range("C12").Value = [double]"2317"
$ws.Range("D12").Value = [double]"1.327251541174455e-07"
$ws.Range("E12").Value = $true
$ws.Range("F12").Value = [double]"-0.5891632373113855"
$ws.Range("G12").Value = [double]"1.30519562593372e-07"
$ws.Range("H12").Value = "Dunn"
$ws.Range("I12").Value = [double]"82.85648148148147"
$ws.Range("J12").Value = "mannwhitneyu"

$ws.Range("A13").Value = "June26"
$ws.Range("B13").Value = "RMS"
$ws.Range("C13").Value = [double]"2175"
$ws.Range("D13").Value = [double]"1.070083676788224e-05"
$ws.Range("E13").Value = $true
$ws.Range("F13").Value = [double]"-0.4917695473251029"
$ws.Range("G13").Value = [double]"1.055027958128037e-05"
$ws.Range("H13").Value = "Dunn"
$ws.Range("I13").Value = [double]"17.32320679012345"
$ws.Range("J13").Value = "mannwhitneyu"

$ws.Range("A14").Value = "June26"
$ws.Range("B14").Value = "tau"
$ws.Range("C14").Value = [double]"619"
$ws.Range("D14").Value = [double]"2.57546540999095e-07"
$ws.Range("E14").Value = $true
$ws.Range("F14").Value = [double]"0.575445816186557"
$ws.Range("G14").Value = [double]"2.533594612822144e-07"
$ws.Range("H14").Value = "Dunn"
$ws.Range("I14").Value = [double]"-42.61480246913578"
$ws.Range("J14").Value = "mannwhitneyu"

$ws.Range("A15").Value = "June26"
$ws.Range("B15").Value = "AUC"
$ws.Range("C15").Value = [double]"1869"
$ws.Range("D15").Value = [double]"0.011659015939608"
$ws.Range("E15").Value = $true
$ws.Range("F15").Value = [double]"-0.2818930041152263"
$ws.Range("G15").Value = [double]"0.01155757319139736"
$ws.Range("H15").Value = "Dunn"
$ws.Range("I15").Value = [double]"4660.227651234569"
$ws.Range("J15").Value = "mannwhitneyu"

$ws.Range("A16").Value = "May20"
$ws.Range("B16").Value = "reactionTime"
$ws.Range("C16").Value = [double]"1113"
$ws.Range("D16").Value = [double]"0.03027169445723122"
$ws.Range("E16").Value = $true
$ws.Range("F16").Value = [double]"0.2366255144032922"
$ws.Range("G16").Value = [double]"0.03003249573877194"
$ws.Range("H16").Value = "Dunn"
$ws.Range("I16").Value = [double]"-0.4185185185185176"
$ws.Range("J16").Value = "mannwhitneyu"

$ws.Range("A17").Value = "May20"
$ws.Range("B17").Value = "peakTime"
$ws.Range("C17").Value = [double]"916.5"
$ws.Range("D17").Value = [double]"0.0006033293714184281"
$ws.Range("E17").Value = $true
$ws.Range("F17").Value = [double]"0.3713991769547325"
$ws.Range("G17").Value = [double]"0.0005963184172090454"
$ws.Range("H17").Value = "Dunn"
$ws.Range("I17").Value = [double]"-1.617901234567899"
$ws.Range("J17").Value = "mannwhitneyu"

$ws.Range("A18").Value = "May20"
$ws.Range("B18").Value = "difference"
$ws.Range("C18").Value = [double]"1190"
$ws.Range("D18").Value = [double]"0.09812993066764997"
$ws.Range("E18").Value = $false
$ws.Range("F18").Value = [double]"0.1838134430727023"
$ws.Range("G18").Value = [double]"0.09750336954318142"
$ws.Range("H18").Value = "Dunn"
$ws.Range("I18").Value = [double]"-1.199382716049378"
$ws.Range("J18").Value = "mannwhitneyu"

$ws.Range("A19").Value = "May20"
$ws.Range("B19").Value = "peakValue"
$ws.Range("C19").Value = [double]"2794.5"
$ws.Range("D19").Value = [double]"2.228994531060003e-16"
$ws.Range("E19").Value = $true
$ws.Range("F19").Value = [double]"-0.9166666666666667"
$ws.Range("G19").Value = [double]"2.172680496951096e-16"
$ws.Range("H19").Value = "Dunn"
$ws.Range("I19").Value = [double]"63.34475308641974"
$ws.Range("J19").Value = "mannwhitneyu"

$ws.Range("A20").Value = "May20"
$ws.Range("B20").Value = "RMS"
$ws.Range("C20").Value = [double]"2526"
$ws.Range("D20").Value = [double]"5.408245471040276e-11"
$ws.Range("E20").Value = $true
$ws.Range("F20").Value = [double]"-0.7325102880658436"
$ws.Range("G20").Value = [double]"5.297948848634982e-11"
$ws.Range("H20").Value = "Dunn"
$ws.Range("I20").Value = [double]"12.87881172839507"
$ws.Range("J20").Value = "mannwhitneyu"

$ws.Range("A21").Value = "May20"
$ws.Range("B21").Value = "tau"
$ws.Range("C21").Value = [double]"608"
$ws.Range("D21").Value = [double]"1.791962740654362e-07"
$ws.Range("E21").Value = $true
$ws.Range("F21").Value = [double]"0.5829903978052127"
$ws.Range("G21").Value = [double]"1.762475032654718e-07"
$ws.Range("H21").Value = "Dunn"
$ws.Range("I21").Value = [double]"-48.5032098765432"
$ws.Range("J21").Value = "mannwhitneyu"

$ws.Range("A22").Value = "May20"
$ws.Range("B22").Value = "AUC"
$ws.Range("C22").Value = [double]"2040"
$ws.Range("D22").Value = [double]"0.0003529112438522091"
$ws.Range("E22").Value = $true
$ws.Range("F22").Value = [double]"-0.3991769547325104"
$ws.Range("G22").Value = [double]"0.0003487918400231757"
$ws.Range("H22").Value = "Dunn"
$ws.Range("I22").Value = [double]"3766.802524691362"
$ws.Range("J22").Value = "mannwhitneyu"

# Remove the now-unused dunn_d column (K), shifting dimension to A1:J22
$ws.Range("K1:K22").Delete()
